# Applies updated cryptocurrency price/volume data as described by the commit
# "Updated cryptos list on Mon Jun 24 21:46:15 UTC 2024 with GitHub Actions".
# All target cells are plain text (Coin name / Link / Price / Volume), so we
# force text number-format before assignment to avoid Excel auto-converting
# numeric-looking strings (e.g. "1.00", "563.40") into floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

Set-TextValue "D2" '59.887.34'
Set-TextValue "E2" '  -6.18%  '

Set-TextValue "D3" '3.345.23'
Set-TextValue "E3" '  -2.49%  '

Set-TextValue "D4" '1.00'
Set-TextValue "E4" '  +0.07%  '

Set-TextValue "D5" '563.40'
Set-TextValue "E5" '  -3.23%  '

Set-TextValue "D6" '130.36'
Set-TextValue "E6" '  -0.21%  '

Set-TextValue "E7" '  +0.05%  '

Set-TextValue "D8" '3.345.17'
Set-TextValue "E8" '  -2.51%  '

Set-TextValue "D9" '0.472'
Set-TextValue "E9" '  -1.82%  '

Set-TextValue "D10" '7.43'
Set-TextValue "E10" '  -2.13%  '

Set-TextValue "E11" '  -6.58%  '

Set-TextValue "E12" '  -2.22%  '

Set-TextValue "D13" '3.913.09'
Set-TextValue "E13" '  -2.37%  '

Set-TextValue "E14" '  -0.19%  '

Set-TextValue "D15" '3.344.55'
Set-TextValue "E15" '  -2.39%  '

Set-TextValue "E16" '  -5.11%  '

Set-TextValue "D17" '24.59'
Set-TextValue "E17" '  -1.57%  '

Set-TextValue "D18" '60.157.30'
Set-TextValue "E18" '  -5.74%  '

Set-TextValue "D19" '5.66'
Set-TextValue "E19" '  -0.36%  '

Set-TextValue "E20" '  +0.63%  '

Set-TextValue "D21" '8.98'
Set-TextValue "E21" '  -9.10%  '

Set-TextValue "D22" '353.72'
Set-TextValue "E22" '  -8.12%  '

Set-TextValue "D23" '0.557'
Set-TextValue "E23" '  -1.61%  '

Set-TextValue "B24" 'WrappedeETH'
Set-TextValue "C24" 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
Set-TextValue "D24" '3.476.78'
Set-TextValue "E24" '  -2.50%  '

Set-TextValue "B25" 'Dai'
Set-TextValue "C25" 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue "D25" '1.00'
Set-TextValue "E25" '  -0.05%  '

Set-TextValue "D26" '69.23'
Set-TextValue "E26" '  -6.32%  '

Set-TextValue "E27" '  +1.77%  '

Set-TextValue "D28" '1.59'
Set-TextValue "E28" '  +11.78%  '

Set-TextValue "D29" '7.44'
Set-TextValue "E29" '  +5.43%  '

Set-TextValue "D30" '0.999'
Set-TextValue "E30" '  +0.42%  '

Set-TextValue "D31" '7.90'
Set-TextValue "E31" '  -1.05%  '

Set-TextValue "E32" '  -4.13%  '

Set-TextValue "D33" '0.153'
Set-TextValue "E33" '  -1.10%  '

Set-TextValue "E34" '  -0.03%  '

Set-TextValue "D35" '3.376.56'
Set-TextValue "E35" '  -2.34%  '

Set-TextValue "D36" '22.88'
Set-TextValue "E36" '  -0.27%  '

Set-TextValue "D37" '5.35'
Set-TextValue "E37" '  +2.82%  '

Set-TextValue "D38" '6.83'
Set-TextValue "E38" '  +0.68%  '

Set-TextValue "E39" '  -0.82%  '

Set-TextValue "D40" '158.07'
Set-TextValue "E40" '  -3.68%  '

Set-TextValue "E41" '  -1.68%  '

Set-TextValue "D42" '0.999'
Set-TextValue "E42" '  -0.01%  '

Set-TextValue "B43" 'ONDO'
Set-TextValue "C43" 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
Set-TextValue "D43" '1.20'
Set-TextValue "E43" '  +7.67%  '

Set-TextValue "B44" 'Filecoin'
Set-TextValue "C44" 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue "D44" '4.37'
Set-TextValue "E44" '  +0.30%  '

Set-TextValue "D45" '40.79'
Set-TextValue "E45" '  -1.72%  '

Set-TextValue "D46" '0.748'
Set-TextValue "E46" '  -5.14%  '

Set-TextValue "D47" '23.54'
Set-TextValue "E47" '  +0.20%  '

Set-TextValue "E48" '  -3.01%  '

Set-TextValue "D49" '6.86'
Set-TextValue "E49" '  +1.73%  '

Set-TextValue "D50" '22.34'
Set-TextValue "E50" '  +9.59%  '

Set-TextValue "D51" '2.41'
Set-TextValue "E51" '  +15.69%  '

